$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new cell at B1, shifting existing B1/C1 to the right (C1/D1).
# xlShiftToRight = -4161
$ws.Range("B1").Insert(-4161)
$ws.Range("B1").Value = '{{xlsxCType t="s"}}{{string}}'

# Apply page setup (paperSize=9/A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
